# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 03:16"

# 2. Update España (row 6) totals
$ws.Range("B6").Value = 21571
$ws.Range("C6").Value = 3494
$ws.Range("E6").Value = 18890

# 3. Update Estados Unidos (row 9) totals
$ws.Range("B9").Value = 19643
$ws.Range("C9").Value = 5854
$ws.Range("E9").Value = 19233

# 4. San Marino overtakes Colombia in the ranking (rows 61/62 swap places)
#    Row 61 becomes San Marino with refreshed totals; row 62 becomes
#    Colombia carrying the figures San Marino previously held at row 61.
$ws.Range("A61").Value = "San Marino"
$ws.Range("B61").Value = 151
$ws.Range("C61").Value = 7
$ws.Range("D61").Value = 4
$ws.Range("E61").Value = 133
$ws.Range("F61").Value = 12
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 14

$ws.Range("A62").Value = "Colombia"
$ws.Range("B62").Value = 145
$ws.Range("C62").Value = 37
$ws.Range("D62").Value = 1
$ws.Range("E62").Value = 144
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 0
